$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1,1).Value = "ISIN"
$ws.Cells.Item(1,2).Value = "Stock Name"
$ws.Cells.Item(1,3).Value = "Mutual Fund"
$ws.Cells.Item(1,4).Value = "Jan_2026"
$ws.Cells.Item(1,5).Value = "Dec_2025"
$ws.Cells.Item(1,6).Value = "Nov_2025"
$ws.Cells.Item(1,7).Value = "MoM"
$ws.Cells.Item(1,8).Value = "QoQ"
$ws.Cells.Item(2,1).Value = "INE931S01010"
$ws.Cells.Item(2,2).Value = "Adani Energy Solutions Limited"
$ws.Cells.Item(2,3).Value = "quant Commodities Fund"
$ws.Cells.Item(2,4).Value = 9.512864
$ws.Cells.Item(2,5).Value = 9.812552999999999
$ws.Cells.Item(2,6).Value = 4.600338
$ws.Cells.Item(2,7).Value = -0.299688999999999
$ws.Cells.Item(2,8).Value = 4.912526000000001
$ws.Cells.Item(3,1).Value = "INE423A01024"
$ws.Cells.Item(3,2).Value = "Adani Enterprises Limited"
$ws.Cells.Item(3,3).Value = "quant Commodities Fund"
$ws.Cells.Item(3,4).Value = 9.024507
$ws.Cells.Item(3,5).Value = 8.987818000000001
$ws.Cells.Item(3,6).Value = 1.884378
$ws.Cells.Item(3,7).Value = 0.03668899999999908
$ws.Cells.Item(3,8).Value = 7.140129
$ws.Cells.Item(4,1).Value = "INE371A01025"
$ws.Cells.Item(4,2).Value = "Graphite India Limited"
$ws.Cells.Item(4,3).Value = "quant Commodities Fund"
$ws.Cells.Item(4,4).Value = 8.991999
$ws.Cells.Item(4,5).Value = 10.366244
$ws.Cells.Item(4,6).Value = 8.969113
$ws.Cells.Item(4,7).Value = -1.374245
$ws.Cells.Item(4,8).Value = 0.02288599999999974
$ws.Cells.Item(5,1).Value = "INE364U01010"
$ws.Cells.Item(5,2).Value = "Adani Green Energy Limited"
$ws.Cells.Item(5,3).Value = "quant Commodities Fund"
$ws.Cells.Item(5,4).Value = 8.815982
$ws.Cells.Item(5,5).Value = 8.948149000000001
$ws.Cells.Item(5,6).Value = 9.243482999999999
$ws.Cells.Item(5,7).Value = -0.1321670000000008
$ws.Cells.Item(5,8).Value = -0.4275009999999995
$ws.Cells.Item(6,1).Value = "INE281B01032"
$ws.Cells.Item(6,2).Value = "Lloyds Metals And Energy Limited"
$ws.Cells.Item(6,3).Value = "quant Commodities Fund"
$ws.Cells.Item(6,4).Value = 8.211117
$ws.Cells.Item(6,5).Value = 8.830628000000001
$ws.Cells.Item(6,6).Value = 5.478523
$ws.Cells.Item(6,7).Value = -0.619511000000001
$ws.Cells.Item(6,8).Value = 2.732594
$ws.Cells.Item(7,1).Value = "INE079A01024"
$ws.Cells.Item(7,2).Value = "Ambuja Cements Ltd"
$ws.Cells.Item(7,3).Value = "quant Commodities Fund"
$ws.Cells.Item(7,4).Value = 5.6717
$ws.Cells.Item(7,5).Value = 5.557017
$ws.Cells.Item(7,6).Value = 0
$ws.Cells.Item(7,7).Value = 0.1146830000000003
$ws.Cells.Item(7,8).Value = 5.6717
$ws.Cells.Item(8,1).Value = "INE331A01037"
$ws.Cells.Item(8,2).Value = "The Ramco Cements Limited"
$ws.Cells.Item(8,3).Value = "quant Commodities Fund"
$ws.Cells.Item(8,4).Value = 4.791299
$ws.Cells.Item(8,5).Value = 0
$ws.Cells.Item(8,6).Value = 0
$ws.Cells.Item(8,7).Value = 4.791299
$ws.Cells.Item(8,8).Value = 4.791299
$ws.Cells.Item(9,1).Value = "INE907A01026"
$ws.Cells.Item(9,2).Value = "Kalyani Steels Ltd"
$ws.Cells.Item(9,3).Value = "quant Commodities Fund"
$ws.Cells.Item(9,4).Value = 4.422661
$ws.Cells.Item(9,5).Value = 4.310524
$ws.Cells.Item(9,6).Value = 4.309386
$ws.Cells.Item(9,7).Value = 0.1121369999999997
$ws.Cells.Item(9,8).Value = 0.1132749999999998
$ws.Cells.Item(10,1).Value = "INE081A01020"
$ws.Cells.Item(10,2).Value = "Tata Steel Limited"
$ws.Cells.Item(10,3).Value = "quant Commodities Fund"
$ws.Cells.Item(10,4).Value = 3.04389
$ws.Cells.Item(10,5).Value = 0
$ws.Cells.Item(10,6).Value = 0
$ws.Cells.Item(10,7).Value = 3.04389
$ws.Cells.Item(10,8).Value = 3.04389
$ws.Cells.Item(11,1).Value = "INE628A01036"
$ws.Cells.Item(11,2).Value = "UPL Limited"
$ws.Cells.Item(11,3).Value = "quant Commodities Fund"
$ws.Cells.Item(11,4).Value = 2.922882
$ws.Cells.Item(11,5).Value = 2.966176
$ws.Cells.Item(11,6).Value = 0
$ws.Cells.Item(11,7).Value = -0.04329399999999994
$ws.Cells.Item(11,8).Value = 2.922882
$ws.Cells.Item(12,1).Value = "INE090A01021"
$ws.Cells.Item(12,2).Value = "ICICI Bank Limited"
$ws.Cells.Item(12,3).Value = "quant Commodities Fund"
$ws.Cells.Item(12,4).Value = 2.717233
$ws.Cells.Item(12,5).Value = 0
$ws.Cells.Item(12,6).Value = 0
$ws.Cells.Item(12,7).Value = 2.717233
$ws.Cells.Item(12,8).Value = 2.717233
$ws.Cells.Item(13,1).Value = "INE814H01029"
$ws.Cells.Item(13,2).Value = "Adani Power Limited"
$ws.Cells.Item(13,3).Value = "quant Commodities Fund"
$ws.Cells.Item(13,4).Value = 0
$ws.Cells.Item(13,5).Value = 2.097552
$ws.Cells.Item(13,6).Value = 2.164396
$ws.Cells.Item(13,7).Value = -2.097552
$ws.Cells.Item(13,8).Value = -2.164396
$ws.Cells.Item(14,1).Value = "INE813H01021"
$ws.Cells.Item(14,2).Value = "Torrent Power Limited"
$ws.Cells.Item(14,3).Value = "quant Commodities Fund"
$ws.Cells.Item(14,4).Value = 0
$ws.Cells.Item(14,5).Value = 1.246209
$ws.Cells.Item(14,6).Value = 1.25413
$ws.Cells.Item(14,7).Value = -1.246209
$ws.Cells.Item(14,8).Value = -1.25413
$ws.Cells.Item(15,1).Value = "INE752E01010"
$ws.Cells.Item(15,2).Value = "Power Grid Corporation of India Limited"
$ws.Cells.Item(15,3).Value = "quant Commodities Fund"
$ws.Cells.Item(15,4).Value = 0
$ws.Cells.Item(15,5).Value = 0
$ws.Cells.Item(15,6).Value = 1.003754
$ws.Cells.Item(15,7).Value = 0
$ws.Cells.Item(15,8).Value = -1.003754
$ws.Cells.Item(16,1).Value = "INE423A20016"
$ws.Cells.Item(16,2).Value = "Adani Enterprises Limited Rights"
$ws.Cells.Item(16,3).Value = "quant Commodities Fund"
$ws.Cells.Item(16,4).Value = 0
$ws.Cells.Item(16,5).Value = 0
$ws.Cells.Item(16,6).Value = 0.043863
$ws.Cells.Item(16,7).Value = 0
$ws.Cells.Item(16,8).Value = -0.043863
$ws.Cells.Item(17,1).Value = "INE075A01022"
$ws.Cells.Item(17,2).Value = "Wipro Ltd"
$ws.Cells.Item(17,3).Value = "quant Commodities Fund"
$ws.Cells.Item(17,4).Value = 0
$ws.Cells.Item(17,5).Value = 1.816168
$ws.Cells.Item(17,6).Value = 0
$ws.Cells.Item(17,7).Value = -1.816168
$ws.Cells.Item(17,8).Value = 0
$ws.Cells.Item(18,1).Value = "INE267A01025"
$ws.Cells.Item(18,2).Value = "Hindustan Zinc Limited"
$ws.Cells.Item(18,3).Value = "quant Commodities Fund"
$ws.Cells.Item(18,4).Value = 0
$ws.Cells.Item(18,5).Value = 6.836897
$ws.Cells.Item(18,6).Value = 0
$ws.Cells.Item(18,7).Value = -6.836897
$ws.Cells.Item(18,8).Value = 0
$ws.Cells.Item(19,1).Value = "INE257A01026"
$ws.Cells.Item(19,2).Value = "Bharat Heavy Electricals Ltd"
$ws.Cells.Item(19,3).Value = "quant Commodities Fund"
$ws.Cells.Item(19,4).Value = 0
$ws.Cells.Item(19,5).Value = 2.903893
$ws.Cells.Item(19,6).Value = 0
$ws.Cells.Item(19,7).Value = -2.903893
$ws.Cells.Item(19,8).Value = 0
$ws.Cells.Item(20,1).Value = "INE245A01021"
$ws.Cells.Item(20,2).Value = "Tata Power Company Limited"
$ws.Cells.Item(20,3).Value = "quant Commodities Fund"
$ws.Cells.Item(20,4).Value = 0
$ws.Cells.Item(20,5).Value = 1.895742
$ws.Cells.Item(20,6).Value = 9.949608
$ws.Cells.Item(20,7).Value = -1.895742
$ws.Cells.Item(20,8).Value = -9.949608
$ws.Cells.Item(21,1).Value = "INE213A01029"
$ws.Cells.Item(21,2).Value = "Oil and Natural Gas Corporation Ltd."
$ws.Cells.Item(21,3).Value = "quant Commodities Fund"
$ws.Cells.Item(21,4).Value = 0
$ws.Cells.Item(21,5).Value = 0
$ws.Cells.Item(21,6).Value = 7.554846
$ws.Cells.Item(21,7).Value = 0
$ws.Cells.Item(21,8).Value = -7.554846
$ws.Cells.Item(22,1).Value = "INE200A01026"
$ws.Cells.Item(22,2).Value = "GE Vernova T&D India Limited"
$ws.Cells.Item(22,3).Value = "quant Commodities Fund"
$ws.Cells.Item(22,4).Value = 0
$ws.Cells.Item(22,5).Value = 2.837839
$ws.Cells.Item(22,6).Value = 0
$ws.Cells.Item(22,7).Value = -2.837839
$ws.Cells.Item(22,8).Value = 0
$ws.Cells.Item(23,1).Value = "INE129A01019"
$ws.Cells.Item(23,2).Value = "GAIL (India) Limited"
$ws.Cells.Item(23,3).Value = "quant Commodities Fund"
$ws.Cells.Item(23,4).Value = 0
$ws.Cells.Item(23,5).Value = 0
$ws.Cells.Item(23,6).Value = 6.487008
$ws.Cells.Item(23,7).Value = 0
$ws.Cells.Item(23,8).Value = -6.487008
$ws.Cells.Item(24,1).Value = "INE296A01032"
$ws.Cells.Item(24,2).Value = "Bajaj Finance Limited"
$ws.Cells.Item(24,3).Value = "quant Commodities Fund"
$ws.Cells.Item(24,4).Value = 0
$ws.Cells.Item(24,5).Value = 0
$ws.Cells.Item(24,6).Value = 9.368145
$ws.Cells.Item(24,7).Value = 0
$ws.Cells.Item(24,8).Value = -9.368145
